$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-11
# from 45212 (2023-10-13) to 45221 (2023-10-22)
for ($row = 2; $row -le 11; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45212) {
        $cell.Value2 = 45221
    }
}
